$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Section header "A: SALARY" -> "SALARY" (drop the "A: " prefix)
$ws.Range("A8").Value = "SALARY"

# 2. Split the old "C: FINAL QUOTE & PROFIT / (LOSS)" header (which lived on
#    A28, above "Overhead allowance:") into its own row: row 27 becomes the
#    (de-prefixed) "FINAL QUOTE & PROFIT / (LOSS)" header, formatted like the
#    other big section headers (bold, 12pt - same look as A28 used to have),
#    and A28 itself is cleared out but keeps its header-row style.
$ws.Range("A27").Font.Bold = $true
$ws.Range("A27").Font.Size = 12
$ws.Range("A27").Value = "FINAL QUOTE  & PROFIT / (LOSS)"
$ws.Rows(27).RowHeight = 15.6

$ws.Range("A28").Value = ""

# 3. Selection left on A28 when the file was saved
$ws.Range("A28").Select()
